$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0.64646998174705406
$ws.Range("B2").Value = 0.64646998174705406
$ws.Range("C2").Value = 0.49996694651932405

# Update row 3 values
$ws.Range("A3").Value = 0.33332129715977532
$ws.Range("B3").Value = 0.33332129715977532
$ws.Range("C3").Value = 0.44285911678262413

# Add new row 4 with values
$ws.Range("A4").Value = 0.33335911440082228
$ws.Range("B4").Value = 0.33335911440082228
$ws.Range("C4").Value = 0.00005348165518017163
